$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data refresh re-sorted rows 3-7 by date (column D, ascending).
# Row contents (columns D, I, J, K, L, M, N, P, Q) are shuffled as follows:
#   new row 3 = old row 7
#   new row 4 = old row 5
#   new row 5 = old row 6
#   new row 6 = old row 3
#   new row 7 = old row 4

# New Row 3 (was row 7)
$ws.Range("D3").Value = 44280
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 30
$ws.Range("K3").Value = 25000
$ws.Range("L3").Value = 25000
$ws.Range("M3").Value = 25000
$ws.Range("N3").Value = "$/caja 18 kilos empedrada"
$ws.Range("P3").Value = 1389
$ws.Range("Q3").Value = 18

# New Row 4 (was row 5)
$ws.Range("D4").Value = 44313
$ws.Range("N4").Value = "$/caja 15 kilos empedrada"

# New Row 5 (was row 6)
$ws.Range("K5").Value = 30000
$ws.Range("L5").Value = 30000
$ws.Range("M5").Value = 30000
$ws.Range("N5").Value = "$/caja 20 kilos empedrada"
$ws.Range("P5").Value = 1500
$ws.Range("Q5").Value = 20

# New Row 6 (was row 3)
$ws.Range("D6").Value = 44315
$ws.Range("I6").Value = "Especial"
$ws.Range("J6").Value = 10

# New Row 7 (was row 4)
$ws.Range("D7").Value = 44315
$ws.Range("J7").Value = 20
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 15000
$ws.Range("N7").Value = "$/caja 15 kilos granel"
$ws.Range("P7").Value = 1000
$ws.Range("Q7").Value = 15
